# Add 5 new rows (11-15) of data to column A, matching the formatting of the
# existing data cells (A2:A10), update the selection to B11, and zoom the
# sheet view to 220%.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing data cell (A10, style s="2")
# down into the new cell range before writing values, so the new cells pick
# up the same vertical-centered / wrap-text style used by the rest of the
# numeric column.
$ws.Range("A10").Copy($ws.Range("A11:A15"))

# New values for column A, rows 11-15.
$ws.Range("A11").Value = 2
$ws.Range("A12").Value = 123
$ws.Range("A13").Value = 4
$ws.Range("A14").Value = 325
$ws.Range("A15").Value = 3

# Update the active selection to B11.
$ws.Range("B11").Select()

# Zoom the sheet view to 220%.
$excel.ActiveWindow.Zoom = 220
